$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebuild column A: int.* urls for registrieren / passwort-vergessen /
# gastzugang / login / kundenportal, each group of 10 "normal" networks
# plus gasnetz-hamburg inserted after the 7th entry -----------------------
$companies = @(
    "avacon-netz.de",
    "bayernwerk-netz.de",
    "sh-netz.com",
    "e-dis-netz.de",
    "energienetze-schaafheim.com",
    "hansegas.com",
    "energienetze-bayern.com",
    "energieversorgung-putzbrunn.de",
    "nordnetz.com",
    "avacon-hochdrucknetz.de"
)
$suffixes = @("registrieren", "passwort-vergessen", "gastzugang", "login", "kundenportal")

# Remove every existing hyperlink up front so stale relationships don't
# linger once the cell text underneath them changes.
$ws.Cells.Hyperlinks.Delete()

$row = 1
foreach ($suf in $suffixes) {
    $idx = 0
    foreach ($c in $companies) {
        if ($idx -eq 7) {
            $gasUrl = "http://int.gasnetz-hamburg.de/" + $suf
            $ws.Cells.Item($row, 1).Value = $gasUrl
            $ws.Hyperlinks.Add($ws.Cells.Item($row, 1), $gasUrl)
            $row = $row + 1
        }
        $ws.Cells.Item($row, 1).Value = "http://int." + $c + "/" + $suf
        $row = $row + 1
        $idx = $idx + 1
    }
}

# --- Cosmetic bits: wider column, scrolled-to-top view with C14 selected -
$ws.Columns.Item(1).ColumnWidth = 60.71
$ws.Range("C14").Select()
